$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap "Enterprises (absolute #)" row and "Enterprises density (per 1000 people)" row
# in the first table (rows 10 and 11) ---
$a10 = $ws.Range("A10").Value2
$b10 = $ws.Range("B10").Value2
$c10 = $ws.Range("C10").Value2
$d10 = $ws.Range("D10").Value2

$a11 = $ws.Range("A11").Value2
$b11 = $ws.Range("B11").Value2
$c11 = $ws.Range("C11").Value2
$d11 = $ws.Range("D11").Value2

$ws.Range("A10").Value = $a11
$ws.Range("B10").Value = $b11
$ws.Range("C10").Value = $c11
$ws.Range("D10").Value = $d11

$ws.Range("A11").Value = $a10
$ws.Range("B11").Value = $b10
$ws.Range("C11").Value = $c10
$ws.Range("D11").Value = $d10

# --- Swap the same two rows in the second table (rows 29 and 30) ---
$a29 = $ws.Range("A29").Value2
$d29 = $ws.Range("D29").Value2

$a30 = $ws.Range("A30").Value2
$d30 = $ws.Range("D30").Value2

$ws.Range("A29").Value = $a30
$ws.Range("D29").Value = $d30

$ws.Range("A30").Value = $a29
$ws.Range("D30").Value = $d29
